$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Objetivos text replaced by the "Docentes responsaveis" text ---
$ws.Range("B10").Value = "3444370 - Rita de Cássia Lacerda Brambilla Rodrigues"
$ws.Range("C10").Value = "3444370 - Rita de Cássia Lacerda Brambilla Rodrigues"

# --- Row 13: gains A13 ("Programa resumido:") and B13/C13 become "Semestral" ---
$ws.Range("A11").Copy($ws.Range("A13"))
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# --- Row 14: "Short syllabus:" + English short-syllabus text ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Proceedings of biochemical processes since the steps of medium preparation and sterilization until the steps of products recuperation and characterization."
$ws.Range("C14").Value = "Proceedings of biochemical processes since the steps of medium preparation and sterilization until the steps of products recuperation and characterization."

# --- Row 15: "Programa:" + date text (reuse existing shared string via Copy to dodge date autodetection) ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))
$ws.Rows.Item(15).RowHeight = 120

# --- Row 16: "Syllabus:" + English syllabus text ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1.Fermentation submerged in bioreactor, involving setting, sterilization, medium addition, cultivation monitoring and samples analysis.`n2.Solid state fermentation followed by recuperation and characterization of the produced enzymes.`n3.Project of enzymes purification by software: the groups receive a problem- mixture and present logical sequences of purification with the respective results and their discussion.`n4.Enzymes characterization in relation to molar mass: calibration of a chromatographic column with known proteins and determination of molar mass of problem-enzyme; determination of problem-enzyme by 280nm absorption and by specific activity.`n5.Use of vegetal and microbial origin amylases for starch processing and ethanol fermentation."
$ws.Range("C16").Value = "1.Fermentation submerged in bioreactor, involving setting, sterilization, medium addition, cultivation monitoring and samples analysis.`n2.Solid state fermentation followed by recuperation and characterization of the produced enzymes.`n3.Project of enzymes purification by software: the groups receive a problem- mixture and present logical sequences of purification with the respective results and their discussion.`n4.Enzymes characterization in relation to molar mass: calibration of a chromatographic column with known proteins and determination of molar mass of problem-enzyme; determination of problem-enzyme by 280nm absorption and by specific activity.`n5.Use of vegetal and microbial origin amylases for starch processing and ethanol fermentation."

# --- Row 17: "Avaliação:" only - clear B17/C17 entirely and drop to default height ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Rows.Item(17).AutoFit()

# --- Row 18: "Método:" + new B18/C18 cells (reuse "3444370..." string via Copy from B10/C10) ---
$ws.Range("A18").Value = "Método:"
$ws.Range("B10").Copy($ws.Range("B18"))
$ws.Range("C10").Copy($ws.Range("C18"))
$ws.Rows.Item(18).RowHeight = 60

# --- Row 19: only A19 label changes to "Critério:"; B19/C19 text unchanged ---
$ws.Range("A19").Value = "Critério:"

# --- Row 20: only A20 label changes to "Norma de recuperação:"; B20/C20 text unchanged ---
$ws.Range("A20").Value = "Norma de recuperação:"

# --- Row 21: only A21 label changes to "Bibliografia:"; B21/C21 text unchanged ---
$ws.Range("A21").Value = "Bibliografia:"
$ws.Rows.Item(21).RowHeight = 120

# --- Row 22: "Requisitos:" only - clear B22/C22 entirely and drop to default height ---
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Rows.Item(22).AutoFit()

# --- Row 23: A23 cleared; new B23/C23 cells reuse the "LOT2013..." string via Copy from B24/C24 ---
$ws.Range("A23").Clear()
$ws.Range("B24").Copy($ws.Range("B23"))
$ws.Range("C24").Copy($ws.Range("C23"))
$ws.Rows.Item(23).RowHeight = 30

# --- Row 24: replace with "LOT2017..." requirement text ---
$ws.Range("B24").Value = "LOT2017 -  Enzimologia  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOT2017 -  Enzimologia  (Requisito fraco)`n"

# --- Row 25 no longer exists: delete it entirely ---
$ws.Rows.Item(25).Delete()
